# Pharma_Society_Report.xlsx update
# - Rename "Sheet1" -> "Report"
# - Update membership-count column (B) values for each society; these are
#   entered as text (leading apostrophe forces Excel to store the
#   numeric-looking value as text rather than a number).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Report"

$ws.Range("B2").Value = "'600"
$ws.Range("B3").Value = "'500"
$ws.Range("B4").Value = "'123"
$ws.Range("B5").Value = "'176"
$ws.Range("B6").Value = "'400"
